$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the definition of "Sequence" (row 19, column B) with the expanded text.
$ws.Range("B19").Value = "a particular section of a course with a specific learning purpose. A course is composed of a series of sequences"

# That row now wraps across more lines, so its height grows from 30 to 45.
$ws.Rows.Item(19).RowHeight = 45

# Add a new term row: "Level" / (no definition) / "等级"
$ws.Range("A20").Value = "Level"
$ws.Range("C20").Value = "等级"
